# Split the final "3," paragraph into "3" and append the 3.a answer text,
# then insert a blank paragraph followed by a new paragraph holding the 3.b answer.
$d = $word.ActiveDocument

# "3," -> "3"
$d.Content.Find.Execute("3,", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2)

# Append the rest of the 3.a sentence right after the "3" in the same paragraph.
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter(".a Ho definito la procedura commentandola e inserendo un interfaccia nel main. Ho cercato di testare l’IsoTriangle e l’IsoTrapezoid ma solo il primo funzionava correttamente.")

# Insert a new, empty paragraph after the 3.a paragraph.
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# Insert another new paragraph (after the blank one) and fill it with the 3.b answer.
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p = $d.Paragraphs.Last
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter("3.b Il test eseguito ha utiliazzato tutti i costruttori di Isotriangle, tutti i getters, setter e gli operatori. Il riscontro è positivo di tutte le parti citate qui sopra e approfondite nel main con i vari commenti.")
